# Committed Corporate Customer excel file
#
# Inserts 4 new columns (C:F) before the existing "T.C (Azure)" column,
# adds the new header labels, switches the Principal cell on row 2 from
# a numeric term value to the text "1M", and leaves the new cells blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing C:E columns ("T.C (Azure)", "T.C (Desc.)", "Error")
# to the right by inserting 4 blank columns at C:F.
$ws.Range("C:F").EntireColumn.Insert()

# Update the Principal value on the data row to a text tenor code.
$ws.Range("B2").Value = "1M"

# Populate the headers for the newly inserted columns.
$ws.Range("C1").Value = "CUST.REMARKS:1"
$ws.Range("D1").Value = "FIQAH"
$ws.Range("E1").Value = "INTEND.DATE"
$ws.Range("F1").Value = "EXP.DATE"

# Match the saved selection/active cell.
$ws.Range("F1").Select()
